# Add a new row documenting "4. Median of Two Sorted Arrays" (Binary Search)
# underneath the existing "981. Time Based Key-Value Store" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 29

# Fill the new row's values (order mirrors the order the strings were
# first typed in: solution summary, link, then question title, then DS)
$ws.Cells.Item($row, 3).Value = "This question is EXTREMELY HARD refer to neetcodes video for explaination"
$ws.Cells.Item($row, 4).Value = "https://youtu.be/q6IEA26hvXc?si=RB1SByCLUSeiK4li "
$ws.Cells.Item($row, 2).Value = "4. Median of Two Sorted Arrays"
$ws.Cells.Item($row, 1).Value = "Binary Search"

# Match the formatting used on the other "Binary Search" rows:
#  - Column A: plain wrapped cell (same as A23/A24/A26/A27)
#  - Column B: bold question text (same as B15/B22/B25)
#  - Column C: green-filled solution summary (same as C22)
#  - Column D: hyperlink style (same as D21/D26)
$ws.Range("A23").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("B15").Copy()
$ws.Range("B29").PasteSpecial(-4122)
$ws.Range("C22").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("D26").Copy()
$ws.Range("D29").PasteSpecial(-4122)

$ws.Rows.Item($row).RowHeight = 28.8

# Turn the URL text in D29 into a real hyperlink, like D21/D26
$ws.Hyperlinks.Add($ws.Range("D29"), "https://youtu.be/q6IEA26hvXc?si=RB1SByCLUSeiK4li ") | Out-Null

# Re-apply the hyperlink cell style (Hyperlinks.Add resets formatting)
$ws.Range("D26").Copy()
$ws.Range("D29").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Range("C31").Select() | Out-Null
